$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "mosajgohar2"
$ws.Range("B3").Value = "wrong password"

$ws.Range("C6").Select()
